# Build site at 2022-09-26 16:07:08 UTC
# The published course-data row ("6310316 - Liana Alvares Rodrigues") was
# re-scraped and the extra value-only row that used to hold it (old row 13,
# with no label in column A) was removed; everything below shifted up one
# row and several of the text values were regenerated in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray row (old row 13: only B13/C13 = "6310316 - Liana Alvares
# Rodrigues", no label) - this shifts rows 14-24 up to 13-23 and carries the
# row heights along correctly.
$ws.Rows("13").Delete()

# --- Fix up the cell values that the new layout places differently ---

# Row 10 (Objetivos:) now shows the docente info instead of the old
# "Introduzir o aluno..." text.
$ws.Range("B10").Value = "6310316 - Liana Alvares Rodrigues"
$ws.Range("C10").Value = "6310316 - Liana Alvares Rodrigues"

# Row 13 (Programa resumido:) now shows "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date "01/01/2012". Use
# Copy/PasteSpecial(values) from a cell that already stores this text so it
# stays a text value instead of being reinterpreted as a date serial.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 (Método:) now shows the docente info.
$ws.Range("B18").Value = "6310316 - Liana Alvares Rodrigues"
$ws.Range("C18").Value = "6310316 - Liana Alvares Rodrigues"

# Row 19 (Critério:) now shows the "Duas provas escritas..." text.
$ws.Range("B19").Value = "Duas provas escritas (P1 e P2) e trabalhos relacionados à disciplina (TRAB)."
$ws.Range("C19").Value = "Duas provas escritas (P1 e P2) e trabalhos relacionados à disciplina (TRAB)."

# Row 20 (Norma de recuperação:) now shows the "Média da Primeira
# Avaliação..." text.
$ws.Range("B20").Value = "Média da Primeira Avaliação = (I)  Prova P1=50%; (II)  Prova P2=50% e (III) `nObs: Fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos."
$ws.Range("C20").Value = "Média da Primeira Avaliação = (I)  Prova P1=50%; (II)  Prova P2=50% e (III) `nObs: Fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos."

# Row 21 (Bibliografia:) now shows the "Será a média aritmética..." text.
$ws.Range("B21").Value = "Será a média aritmética da nota do aluno na primeira avaliação e da nota do aluno numa prova escrita na recuperação"
$ws.Range("C21").Value = "Será a média aritmética da nota do aluno na primeira avaliação e da nota do aluno numa prova escrita na recuperação"
